$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add new column headers I1 ("I0") and J1 ("IF") ---
# Copy formatting from the existing header cell H1 so the new headers
# share the same style (bold font, border, centered alignment).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-78: add values for the new I and J columns ---
$rows = @(
    @(2,8,8),
    @(3,8,8),
    @(4,8,8),
    @(5,8,8),
    @(6,7,8),
    @(7,8,8),
    @(8,8,8),
    @(9,8,8),
    @(10,8,9),
    @(11,8,8),
    @(12,7,7),
    @(13,8,8),
    @(14,7,7),
    @(15,9,9),
    @(16,7,7),
    @(17,8,8),
    @(18,9,9),
    @(19,9,9),
    @(20,9,9),
    @(21,8,9),
    @(22,9,9),
    @(23,7,7),
    @(24,7,7),
    @(25,8,8),
    @(26,7,8),
    @(27,8,8),
    @(28,9,9),
    @(29,8,8),
    @(30,9,9),
    @(31,8,8),
    @(32,8,8),
    @(33,8,8),
    @(34,7,8),
    @(35,9,9),
    @(36,8,8),
    @(37,8,8),
    @(38,8,8),
    @(39,7,8),
    @(40,8,8),
    @(41,8,9),
    @(42,8,8),
    @(43,7,7),
    @(44,9,9),
    @(45,8,8),
    @(46,7,8),
    @(47,8,8),
    @(48,7,8),
    @(49,7,7),
    @(50,8,8),
    @(51,7,8),
    @(52,8,8),
    @(53,8,8),
    @(54,7,7),
    @(55,8,8),
    @(56,7,7),
    @(57,7,8),
    @(58,9,9),
    @(59,7,8),
    @(60,8,8),
    @(61,8,9),
    @(62,8,8),
    @(63,9,9),
    @(64,9,9),
    @(65,10,10),
    @(66,7,7),
    @(67,7,7),
    @(68,5,6),
    @(69,6,7),
    @(70,6,6),
    @(71,6,6),
    @(72,4,4),
    @(73,5,6),
    @(74,4,4),
    @(75,8,8),
    @(76,7,7),
    @(77,4,4),
    @(78,4,4)
)

foreach ($row in $rows) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
